$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell values (row 1 title, row 2 headers)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Teki Susu - de 12 de Agosto de 2024"
$ws.Range("B1").Value = ""
$ws.Range("C1").Value = ""
$ws.Range("F1").Value = ""
$ws.Range("G1").Value = ""

$ws.Range("A2").Value = "Nombre"
$ws.Range("B2").Value = "Semana"
$ws.Range("C2").Value = "Lunes"
$ws.Range("D2").Value = "Martes"
$ws.Range("E2").Value = "Miercoles"
$ws.Range("F2").Value = "Jueves"
$ws.Range("G2").Value = "Viernes"
$ws.Range("H2").Value = "Sábado"
$ws.Range("I2").Value = "Estado"

# ---------------------------------------------------------------------------
# 2. Formatting - reuse existing styles via copy/paste-special (formats only)
#    rather than setting ad-hoc formatting properties, so we don't create
#    brand-new style-table entries.
# ---------------------------------------------------------------------------

# Title-row style (currently on A1:C1) -> stretch across D1:I1 too.
$ws.Range("A1").Copy()
$ws.Range("D1:I1").PasteSpecial(-4122)

# Header-row style (currently on A2:G2) -> extend to the two new columns.
$ws.Range("A2").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 16
$ws.Columns.Item(8).ColumnWidth = 14
$ws.Columns.Item(9).ColumnWidth = 20.42578125

# ---------------------------------------------------------------------------
# 4. Merge title row across the full width
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Merge()

# ---------------------------------------------------------------------------
# 5. AutoFilter over the new header range
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A2:I2").AutoFilter()

# ---------------------------------------------------------------------------
# 6. Defined name _FilterDatabase -> new range
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Hoja1!`$A`$2:`$I`$2"
    }
}

# ---------------------------------------------------------------------------
# 7. Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("K8").Select()
$excel.ActiveWindow.WindowState = -4143
$wb.Windows.Item(1).WindowWidth = 20490
$wb.Windows.Item(1).WindowHeight = 7650
